$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.203.30"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.877.10"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'316.53"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4317"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D8").Value = "'0.3704"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "'0.07424"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").Value = "'0.8871"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'21.21"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "1.917.96"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "'5.487"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "'6.642"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'0.06994"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'81.22"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "'0.000009140"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'15.64"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "28.158.19"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'5.097"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").Value = "2.151.88"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "'1.985"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "'154.43"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'18.70"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "'5.456"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").Value = "'118.32"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "'1.905"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "'0.08991"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'0.7988"
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("D33").Value = "'4.699"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "'1.176"
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("D35").Value = "'2.982"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'1.137"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'0.05481"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'2.894"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").Value = "'0.1700"
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("D42").Value = "'0.5180"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'6.876"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "'8.599"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "'10.58"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").Value = "'0.06606"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "'0.4777"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'105.70"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "'1.003"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "'1.662"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "'1.847"
$ws.Range("E51").Value = "  +5.18%  "
